$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the order rows: 10 "Jean" orders, one per year from 2021-2030,
# with incrementing delivery time, zeroed quantity columns, a flat 750
# final value, an "asdf" message and "Concluído" status.
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $year = 2021 + $i
    $hora = "00:{0:D2}" -f $i

    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 2).Value = "Jean"
    $ws.Cells.Item($r, 3).Value = "30/12/$year"
    $ws.Cells.Item($r, 4).Value = $hora
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 750
    $ws.Cells.Item($r, 10).Value = "asdf"
    $ws.Cells.Item($r, 11).Value = "Concluído"
}

# Column C only needs to fit the "dd/mm/yyyy" strings now (was sized for
# "Nome cliente"-length text before).
$ws.Columns("C").ColumnWidth = 10.5

# Zoom in a bit and move the live selection like the author left it.
$excel.ActiveWindow.Zoom = 130
[void]$ws.Range("M9").Select()
